# updates with supplemental figures
# Refresh the Fixed-Residue PSSM score matrix (B2:K21) with the recomputed
# values. The "background"/unused-residue placeholder score moves from
# -19.83129332503361 to -15.23717260797651, and every populated score cell
# is replaced with its newly recalculated value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-15.23717260797651,1.24558808988386,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651),
    @(-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,1.60489531752938,-15.23717260797651,-15.23717260797651),
    @(-15.23717260797651,1.898947386999727,-15.23717260797651,-15.23717260797651,3.676572910641023,-15.23717260797651,1.84340553295301,-15.23717260797651,2.599598963626845,-15.23717260797651),
    @(-15.23717260797651,1.529006868543322,-15.23717260797651,-15.23717260797651,-15.23717260797651,3.334391184860558,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651),
    @(-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651),
    @(2.973693707129837,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651),
    @(-15.23717260797651,-15.23717260797651,-15.23717260797651,2.095816724995375,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651),
    @(3.602181904150547,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651),
    @(-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,0.792567547226802,-15.23717260797651,2.093192863172357),
    @(-15.23717260797651,-15.23717260797651,-15.23717260797651,2.525112095369227,-15.23717260797651,2.04711833382781,-15.23717260797651,-15.23717260797651,-15.23717260797651,1.501038428364564),
    @(-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651),
    @(-15.23717260797651,-15.23717260797651,-15.23717260797651,1.925543243453868,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,1.48137364613752,2.379383414800549),
    @(-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,1.969507928308031),
    @(-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651),
    @(-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,1.917259860814394,-15.23717260797651),
    @(-15.23717260797651,2.431968513891992,4.321892608851974,-15.23717260797651,-15.23717260797651,-15.23717260797651,1.962815981800988,0.9184559983760019,2.085292523159054,-15.23717260797651),
    @(-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,1.979393280972999,0.6122569286559247,1.64414452814219,-15.23717260797651),
    @(-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,-15.23717260797651,1.979582721662504,1.649393257011339,-15.23717260797651,-15.23717260797651),
    @(-15.23717260797651,1.772579482255223,-15.23717260797651,-15.23717260797651,2.850564731935308,-15.23717260797651,1.194276026348321,3.116030227028791,-15.23717260797651,1.919433919277306),
    @(-15.23717260797651,1.137680092001649,-15.23717260797651,2.625252872088421,-15.23717260797651,2.531069880476027,1.224962170242821,-15.23717260797651,-15.23717260797651,-15.23717260797651)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $colNum = $j + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}

Write-Host "Updated PSSM values"